# Fruta / hortaliza, semanal
# Insert a new weekly record as row 208, pushing the existing rows
# 208:221 down to 209:222, then populate the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 208 - shifts rows 208:221 down to 209:222
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new weekly observation
$ws.Cells.Item(208, 1).Value2 = 7
$ws.Cells.Item(208, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(208, 3).Value2 = "Ñuble"
$ws.Cells.Item(208, 4).Value2 = 44585
$ws.Cells.Item(208, 5).Value2 = 16
$ws.Cells.Item(208, 6).Value2 = 100112008
$ws.Cells.Item(208, 7).Value2 = "Coliflor"
$ws.Cells.Item(208, 8).Value2 = "Sin especificar"
$ws.Cells.Item(208, 9).Value2 = "Primera"
$ws.Cells.Item(208, 10).Value2 = 300
$ws.Cells.Item(208, 11).Value2 = 850
$ws.Cells.Item(208, 12).Value2 = 900
$ws.Cells.Item(208, 13).Value2 = 875
$ws.Cells.Item(208, 14).Value2 = "$/unidad"
$ws.Cells.Item(208, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(208, 16).Value2 = 875
$ws.Cells.Item(208, 17).Value2 = 1
$ws.Cells.Item(208, 18).Value2 = "Hortaliza"
